# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 872.4375
$ws.Range("I28").Value = 329.91666
$ws.Range("J28").Value = 2500
$ws.Range("K28").Value = 329.91666
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 155.08334
$ws.Range("N28").Value = -3470
$ws.Range("H33").Value = 51542.15
$ws.Range("I33").Value = 67122.734
$ws.Range("K33").Value = 67122.734
$ws.Range("M33").Value = -66893.734
$ws.Range("H86").Value = 4370.5293
$ws.Range("I86").Value = 2140.4
$ws.Range("J86").Value = 7556.4287
$ws.Range("K86").Value = 2140.4
$ws.Range("L86").Value = 7556.4287
$ws.Range("M86").Value = -1017.4
$ws.Range("N86").Value = -9802.4287
$ws.Range("H89").Value = 4370.5293
$ws.Range("I89").Value = 2140.4
$ws.Range("J89").Value = 7556.4287
$ws.Range("K89").Value = 10702
$ws.Range("L89").Value = 37782.14350000001
$ws.Range("M89").Value = -5086
$ws.Range("N89").Value = -49014.14350000001
$ws.Range("H92").Value = 1223.4584
$ws.Range("I92").Value = 1280.5294
$ws.Range("J92").Value = 1084.8572
$ws.Range("K92").Value = 1280.5294
$ws.Range("L92").Value = 1084.8572
$ws.Range("M92").Value = -32.5293999999999
$ws.Range("N92").Value = -3580.8572
$ws.Range("H107").Value = 452.88235
$ws.Range("I107").Value = 399.92307
$ws.Range("J107").Value = 625
$ws.Range("K107").Value = 399.92307
$ws.Range("L107").Value = 625
$ws.Range("M107").Value = 1520.07693
$ws.Range("N107").Value = -4465
$ws.Range("H113").Value = 64425.75
$ws.Range("J113").Value = 1901.4166
$ws.Range("L113").Value = 1901.4166
$ws.Range("N113").Value = -8409.4166
$ws.Range("H118").Value = 7163.5293
$ws.Range("I118").Value = 8391.429
$ws.Range("J118").Value = 1433.3334
$ws.Range("K118").Value = 25174.287
$ws.Range("L118").Value = 4300.0002
$ws.Range("M118").Value = -23517.287
$ws.Range("N118").Value = -7614.0002
$ws.Range("H127").Value = 23811162
$ws.Range("I127").Value = 434
$ws.Range("J127").Value = 33335454
$ws.Range("K127").Value = 1302
$ws.Range("L127").Value = 100006362
$ws.Range("M127").Value = 3658
$ws.Range("N127").Value = -100016282
$ws.Range("H128").Value = 34499.8
$ws.Range("J128").Value = 35374.75
$ws.Range("L128").Value = 35374.75
$ws.Range("N128").Value = -45334.75
$ws.Range("H132").Value = 10008254
$ws.Range("I132").Value = 13900031
$ws.Range("J132").Value = 828.5714
$ws.Range("K132").Value = 41700093
$ws.Range("L132").Value = 2485.7142
$ws.Range("M132").Value = -41697563
$ws.Range("N132").Value = -7545.7142
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140
$ws.Range("H135").Value = 1081.7084
$ws.Range("I135").Value = 753.0222
$ws.Range("K135").Value = 6777.1998
$ws.Range("M135").Value = -4242.1998
$ws.Range("H137").Value = 1060.2034
$ws.Range("I137").Value = 1020.7308
$ws.Range("K137").Value = 3062.1924
$ws.Range("M137").Value = -512.1923999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1047.03
$ws.Range("I32").Value = 1059.5212
$ws.Range("J32").Value = 851.3333
$ws.Range("K32").Value = 1059.5212
$ws.Range("L32").Value = 851.3333
$ws.Range("M32").Value = -772.5211999999999
$ws.Range("N32").Value = -1425.3333
$ws.Range("H45").Value = 47316.273
$ws.Range("I45").Value = 57054.61
$ws.Range("J45").Value = 3493.75
$ws.Range("K45").Value = 57054.61
$ws.Range("L45").Value = 3493.75
$ws.Range("M45").Value = -56677.61
$ws.Range("N45").Value = -4247.75
$ws.Range("H74").Value = 486.18518
$ws.Range("I74").Value = 471.79166
$ws.Range("J74").Value = 601.3333
$ws.Range("K74").Value = 471.79166
$ws.Range("L74").Value = 601.3333
$ws.Range("M74").Value = 402.20834
$ws.Range("N74").Value = -2349.3333
$ws.Range("H77").Value = 486.18518
$ws.Range("I77").Value = 471.79166
$ws.Range("J77").Value = 601.3333
$ws.Range("K77").Value = 2358.9583
$ws.Range("L77").Value = 3006.6665
$ws.Range("M77").Value = 2009.0417
$ws.Range("N77").Value = -11742.6665
$ws.Range("H122").Value = 1113.5454
$ws.Range("I122").Value = 1689.1666
$ws.Range("J122").Value = 422.8
$ws.Range("K122").Value = 5067.4998
$ws.Range("L122").Value = 1268.4
$ws.Range("M122").Value = -2617.4998
$ws.Range("N122").Value = -6168.4
$ws.Range("H132").Value = 3833.4048
$ws.Range("I132").Value = 3710.7104
$ws.Range("K132").Value = 11132.1312
$ws.Range("M132").Value = -8602.1312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 32243.625
$ws.Range("J81").Value = 32243.625
$ws.Range("L81").Value = 32243.625
$ws.Range("N81").Value = -34365.625
$ws.Range("H84").Value = 32243.625
$ws.Range("J84").Value = 32243.625
$ws.Range("L84").Value = 96730.875
$ws.Range("N84").Value = -107338.875
$ws.Range("H134").Value = 2169.0364
$ws.Range("I134").Value = 1873.5471
$ws.Range("J134").Value = 9999.5
$ws.Range("K134").Value = 5620.6413
$ws.Range("L134").Value = 29998.5
$ws.Range("M134").Value = -3085.6413
$ws.Range("N134").Value = -35068.5
$ws.Range("H135").Value = 50874.75
$ws.Range("J135").Value = 50874.75
$ws.Range("L135").Value = 50874.75
$ws.Range("N135").Value = -61014.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1170.1818
$ws.Range("I16").Value = 819.6667
$ws.Range("J16").Value = 2747.5
$ws.Range("K16").Value = 819.6667
$ws.Range("L16").Value = 2747.5
$ws.Range("M16").Value = -532.6667
$ws.Range("N16").Value = -3321.5
$ws.Range("H31").Value = 38253.727
$ws.Range("I31").Value = 2218.1765
$ws.Range("J31").Value = 64888.695
$ws.Range("K31").Value = 2218.1765
$ws.Range("L31").Value = 64888.695
$ws.Range("M31").Value = -1923.1765
$ws.Range("N31").Value = -65478.695
$ws.Range("H34").Value = 38253.727
$ws.Range("I34").Value = 2218.1765
$ws.Range("J34").Value = 64888.695
$ws.Range("K34").Value = 2218.1765
$ws.Range("L34").Value = 64888.695
$ws.Range("M34").Value = -2016.1765
$ws.Range("N34").Value = -65292.695
$ws.Range("H113").Value = 1170.1818
$ws.Range("I113").Value = 819.6667
$ws.Range("J113").Value = 2747.5
$ws.Range("K113").Value = 819.6667
$ws.Range("L113").Value = 2747.5
$ws.Range("M113").Value = 1350.3333
$ws.Range("N113").Value = -7087.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1074.8928
$ws.Range("I5").Value = 787.6429000000001
$ws.Range("J5").Value = 1362.1428
$ws.Range("K5").Value = 2362.9287
$ws.Range("L5").Value = 4086.4284
$ws.Range("M5").Value = -2250.9287
$ws.Range("N5").Value = -4310.428400000001
$ws.Range("H112").Value = 73372.07000000001
$ws.Range("I112").Value = 250550
$ws.Range("J112").Value = 2500.9
$ws.Range("K112").Value = 751650
$ws.Range("L112").Value = 7502.700000000001
$ws.Range("M112").Value = -750542
$ws.Range("N112").Value = -9718.700000000001
$ws.Range("H134").Value = 3637
$ws.Range("I134").Value = 2067.8572
$ws.Range("J134").Value = 7298.3335
$ws.Range("K134").Value = 6203.571599999999
$ws.Range("L134").Value = 21895.0005
$ws.Range("M134").Value = -1133.571599999999
$ws.Range("N134").Value = -32035.0005
$ws.Range("H135").Value = 1074.8928
$ws.Range("I135").Value = 787.6429000000001
$ws.Range("J135").Value = 1362.1428
$ws.Range("K135").Value = 7088.7861
$ws.Range("L135").Value = 12259.2852
$ws.Range("M135").Value = -4553.7861
$ws.Range("N135").Value = -17329.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3389.9167
$ws.Range("J43").Value = 7003.8
$ws.Range("L43").Value = 7003.8
$ws.Range("N43").Value = -7305.8
$ws.Range("H57").Value = 14055
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 43524908
$ws.Range("J80").Value = 2159
$ws.Range("L80").Value = 2159
$ws.Range("N80").Value = -4155
$ws.Range("H83").Value = 43524908
$ws.Range("J83").Value = 2159
$ws.Range("L83").Value = 10795
$ws.Range("N83").Value = -20779
$ws.Range("H102").Value = 2388.64
$ws.Range("I102").Value = 1870.5
$ws.Range("J102").Value = 2866.923
$ws.Range("K102").Value = 1870.5
$ws.Range("L102").Value = 2866.923
$ws.Range("M102").Value = -248.5
$ws.Range("N102").Value = -6110.923
$ws.Range("H113").Value = 1646.1818
$ws.Range("I113").Value = 1866.625
$ws.Range("J113").Value = 1520.2142
$ws.Range("K113").Value = 1866.625
$ws.Range("L113").Value = 1520.2142
$ws.Range("M113").Value = 303.375
$ws.Range("N113").Value = -5860.2142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2870.15
$ws.Range("J7").Value = 4044.4443
$ws.Range("L7").Value = 4044.4443
$ws.Range("N7").Value = -4268.4443
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H126").Value = 2870.15
$ws.Range("J126").Value = 4044.4443
$ws.Range("L126").Value = 12133.3329
$ws.Range("N126").Value = -17073.3329
$ws.Range("H132").Value = 2729.9783
$ws.Range("I132").Value = 2732.8057
$ws.Range("J132").Value = 2719.8
$ws.Range("K132").Value = 8198.417099999999
$ws.Range("L132").Value = 8159.400000000001
$ws.Range("M132").Value = -5668.417099999999
$ws.Range("N132").Value = -13219.4
$ws.Range("H136").Value = 1145.4717
$ws.Range("I136").Value = 926.8
$ws.Range("J136").Value = 2375.5
$ws.Range("K136").Value = 2780.4
$ws.Range("L136").Value = 7126.5
$ws.Range("M136").Value = -230.3999999999996
$ws.Range("N136").Value = -12226.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 333767.34
$ws.Range("I100").Value = 500151
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1000302
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -999761
$ws.Range("N100").Value = -3082
$ws.Range("H132").Value = 2313.75
$ws.Range("I132").Value = 2276.907
$ws.Range("K132").Value = 6830.721
$ws.Range("M132").Value = -4300.721
$ws.Range("H136").Value = 784.0526
$ws.Range("I136").Value = 519.8333
$ws.Range("J136").Value = 1774.875
$ws.Range("K136").Value = 1559.4999
$ws.Range("L136").Value = 5324.625
$ws.Range("M136").Value = 990.5001
$ws.Range("N136").Value = -10424.625
